# Auto-generated edit script: updates market-price-derived columns (H-N)
# across multiple worksheets to match the refreshed source data.
$wb = $excel.ActiveWorkbook


# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1315.8823
$ws.Range("J112").Value = 1340.6061
$ws.Range("L112").Value = 4021.8183
$ws.Range("N112").Value = -6237.8183
$ws.Range("H125").Value = 7594.6665
$ws.Range("I125").Value = 10032
$ws.Range("J125").Value = 7107.2
$ws.Range("K125").Value = 90288
$ws.Range("L125").Value = 63964.8
$ws.Range("M125").Value = -87828
$ws.Range("N125").Value = -68884.79999999999
$ws.Range("H131").Value = 4105.352
$ws.Range("I131").Value = 582.8570999999999
$ws.Range("J131").Value = 4629.9785
$ws.Range("K131").Value = 1748.5713
$ws.Range("L131").Value = 13889.9355
$ws.Range("M131").Value = 3291.4287
$ws.Range("N131").Value = -23969.9355
$ws.Range("H141").Value = 2358
$ws.Range("I141").Value = 1760.25
$ws.Range("J141").Value = 4151.25
$ws.Range("K141").Value = 5280.75
$ws.Range("L141").Value = 12453.75
$ws.Range("M141").Value = -100.75
$ws.Range("N141").Value = -22813.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 175.5
$ws.Range("I4").Value = 201
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 201
$ws.Range("L4").Value = 150
$ws.Range("M4").Value = -85
$ws.Range("N4").Value = -382
$ws.Range("H5").Value = 388.3
$ws.Range("I5").Value = 380.16666
$ws.Range("K5").Value = 380.16666
$ws.Range("M5").Value = -268.16666
$ws.Range("H104").Value = 35893.75
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 35893.75
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 35893.75
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -42881.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 388.3
$ws.Range("I4").Value = 380.16666
$ws.Range("K4").Value = 380.16666
$ws.Range("M4").Value = -265.16666
$ws.Range("H22").Value = 270.4
$ws.Range("I22").Value = 212.5
$ws.Range("J22").Value = 502
$ws.Range("K22").Value = 212.5
$ws.Range("L22").Value = 502
$ws.Range("M22").Value = -39.5
$ws.Range("N22").Value = -848

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 322.47058
$ws.Range("I22").Value = 365.77777
$ws.Range("J22").Value = 273.75
$ws.Range("K22").Value = 365.77777
$ws.Range("L22").Value = 273.75
$ws.Range("M22").Value = -15.77776999999998
$ws.Range("N22").Value = -973.75
$ws.Range("H62").Value = 2525
$ws.Range("I62").Value = 2300
$ws.Range("J62").Value = 2570
$ws.Range("K62").Value = 2300
$ws.Range("L62").Value = 2570
$ws.Range("M62").Value = -1676
$ws.Range("N62").Value = -3818
$ws.Range("H65").Value = 2525
$ws.Range("I65").Value = 2300
$ws.Range("J65").Value = 2570
$ws.Range("K65").Value = 11500
$ws.Range("L65").Value = 12850
$ws.Range("M65").Value = -8380
$ws.Range("N65").Value = -19090
$ws.Range("H99").Value = 2734.28
$ws.Range("I99").Value = 2061.5
$ws.Range("J99").Value = 3050.8823
$ws.Range("K99").Value = 2061.5
$ws.Range("L99").Value = 3050.8823
$ws.Range("M99").Value = -563.5
$ws.Range("N99").Value = -6046.8823
$ws.Range("H107").Value = 1043
$ws.Range("I107").Value = 1148.8667
$ws.Range("J107").Value = 816.1429000000001
$ws.Range("K107").Value = 1148.8667
$ws.Range("L107").Value = 816.1429000000001
$ws.Range("M107").Value = 771.1333
$ws.Range("N107").Value = -4656.1429
$ws.Range("H126").Value = 2734.28
$ws.Range("I126").Value = 2061.5
$ws.Range("J126").Value = 3050.8823
$ws.Range("K126").Value = 6184.5
$ws.Range("L126").Value = 9152.6469
$ws.Range("M126").Value = -3714.5
$ws.Range("N126").Value = -14092.6469

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 19800.31
$ws.Range("I131").Value = 448.25
$ws.Range("J131").Value = 24848.674
$ws.Range("K131").Value = 1344.75
$ws.Range("L131").Value = 74546.022
$ws.Range("M131").Value = 3695.25
$ws.Range("N131").Value = -84626.022

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 117.09524
$ws.Range("I2").Value = 63.76923
$ws.Range("J2").Value = 203.75
$ws.Range("K2").Value = 63.76923
$ws.Range("L2").Value = 203.75
$ws.Range("M2").Value = 49.23077
$ws.Range("N2").Value = -429.75
$ws.Range("H93").Value = 35000
$ws.Range("J93").Value = 35000
$ws.Range("L93").Value = 35000
$ws.Range("N93").Value = -38744
$ws.Range("H107").Value = 777437.75
$ws.Range("I107").Value = 331.33334
$ws.Range("J107").Value = 2525927.2
$ws.Range("K107").Value = 331.33334
$ws.Range("L107").Value = 2525927.2
$ws.Range("M107").Value = 1588.66666
$ws.Range("N107").Value = -2529767.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2761.4375
$ws.Range("I7").Value = 1743.2222
$ws.Range("J7").Value = 4070.5715
$ws.Range("K7").Value = 1743.2222
$ws.Range("L7").Value = 4070.5715
$ws.Range("M7").Value = -1631.2222
$ws.Range("N7").Value = -4294.5715
$ws.Range("H22").Value = 1798.75
$ws.Range("J22").Value = 698.3333
$ws.Range("L22").Value = 698.3333
$ws.Range("N22").Value = -1288.3333
$ws.Range("H27").Value = 1798.75
$ws.Range("J27").Value = 698.3333
$ws.Range("L27").Value = 698.3333
$ws.Range("N27").Value = -912.3333
$ws.Range("H40").Value = 41016.04
$ws.Range("I40").Value = 73645.86
$ws.Range("J40").Value = 2947.9167
$ws.Range("K40").Value = 73645.86
$ws.Range("L40").Value = 2947.9167
$ws.Range("M40").Value = -73509.86
$ws.Range("N40").Value = -3219.9167
$ws.Range("H61").Value = 2125.4666
$ws.Range("I61").Value = 2161
$ws.Range("J61").Value = 2027.75
$ws.Range("K61").Value = 2161
$ws.Range("L61").Value = 2027.75
$ws.Range("M61").Value = -1959
$ws.Range("N61").Value = -2431.75
$ws.Range("H113").Value = 2125.4666
$ws.Range("I113").Value = 2161
$ws.Range("J113").Value = 2027.75
$ws.Range("K113").Value = 2161
$ws.Range("L113").Value = 2027.75
$ws.Range("M113").Value = 9
$ws.Range("N113").Value = -6367.75
$ws.Range("H122").Value = 3139.182
$ws.Range("I122").Value = 3106.5715
$ws.Range("J122").Value = 3196.25
$ws.Range("K122").Value = 9319.7145
$ws.Range("L122").Value = 9588.75
$ws.Range("M122").Value = -6869.7145
$ws.Range("N122").Value = -14488.75
$ws.Range("H126").Value = 2761.4375
$ws.Range("I126").Value = 1743.2222
$ws.Range("J126").Value = 4070.5715
$ws.Range("K126").Value = 5229.6666
$ws.Range("L126").Value = 12211.7145
$ws.Range("M126").Value = -2759.6666
$ws.Range("N126").Value = -17151.7145
$ws.Range("H132").Value = 3359.12
$ws.Range("I132").Value = 3359.12
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10077.36
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7547.360000000001
$ws.Range("N132").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 71881.42999999999
$ws.Range("I107").Value = 449.0909
$ws.Range("J107").Value = 333800
$ws.Range("K107").Value = 1347.2727
$ws.Range("L107").Value = 1001400
$ws.Range("M107").Value = 572.7273
$ws.Range("N107").Value = -1005240
